# Apply the schedule-template updates described in the commit:
#  - Add an "Instructor titular" label/value pair in D5/E5
#  - Insert a new competency row ("Instalación + Manual de Usuario") right
#    after the "Proyecto 2 + Pruebas (Tc.)" row (row 28), pushing the
#    following "Emprendimiento" rows (previously 29-32) down to 30-33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: record the lead instructor for the schedule ---
$ws.Range("D5").Value = "Instructor titular:"
$ws.Range("D5").Font.Bold = $true
$ws.Range("E5").Value = "DANIEL DAVID BENAVIDES SÁNCHEZ "

# --- Insert a new row before row 29 (shifts old rows 29-32 to 30-33) ---
$ws.Rows.Item(29).Insert()

# --- Populate the newly inserted row 29 with the "Instalación + Manual
#     de Usuario" competency details (mirrors the pattern used for the
#     "Proyecto 2 + Pruebas (Tc.)" row right above it) ---
$ws.Range("A29").Value = "Instalación + Manual de Usuario"
$ws.Range("B29").Value = "Desarrollar la solución de software de acuerdo con el diseño y metodologías de desarrollo."
$ws.Range("C29").Value = "Desarrollo de la solución de software."
$ws.Range("D29").Value = "Codificar el software empleando el lenguaje de programación seleccionado."
